# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310, *_new -> *_FV2404
# Then turn the used range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into a proper Table (adds xl/tables/table1.xml,
# the autofilter and tableParts reference) with the just-renamed headers.
$usedRange = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $usedRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split after row 1).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
